$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# "Wario - Level 6 begin"
#
# A new timing row ("Leave Level 5") is inserted right before the existing
# "Level 6" section header (old row 50), pushing everything from the old
# row 50 onward down by one row. The newly freed rows 45-49 (Level 5 pipe /
# key entries) get their start ("B") times filled in for the first time,
# and the "Enter Level 6" row (old row 51, now row 52) also gets its start
# time filled in.
# ---------------------------------------------------------------------------

# 1) Insert a blank row at 50 (this shifts the "Level 6" header and every
#    row below it down by one, and Excel auto-adjusts formulas, merged
#    cells, the used-range dimension, etc.)
$ws.Rows.Item(50).Insert()

# 2) Clone the formatting of the row right above (row 49, "Use Key") into
#    the freshly inserted row 50 so the new row matches the sheet's normal
#    data-row styling.
$ws.Range("A49:D49").Copy()
$ws.Range("A50:D50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the previously-blank "B" (start time) values for the last few
#    Level 5 entries.
$ws.Range("B45").Value() = 15841
$ws.Range("B46").Value() = 16245
$ws.Range("B47").Value() = 16760
$ws.Range("B48").Value() = 16939
$ws.Range("B49").Value() = 17105

# 4) Populate the new row with the "Leave Level 5" entry.
$ws.Range("A50").Value() = "Leave Level 5"
$ws.Range("B50").Value() = 17309
$ws.Range("C50").Value() = 18801
$ws.Range("D50").Formula = "=IF(B50>0,C50-B50,0)"

# 5) The "Enter Level 6" row (now row 52) also gets its start time filled
#    in for the first time.
$ws.Range("B52").Value() = 17558

# 6) Leave the selection where the author ended up after the edit.
$ws.Activate()
$ws.Range("B53").Select()
